$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data blocks of rows 2-4 ("Critical low battery ...") and
# rows 5-7 ("Compass Error ...") so that the Compass Error block now
# appears first (rows 2-4) and the Critical low battery block appears
# second (rows 5-7). Columns A, E, F are left untouched.

$compassMessage = "Compass Error Compass Error Compass disconnected ."
$criticalMessage = "Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress ."

# New rows 2-4: Compass Error block
$ws.Range("B2").Value = $compassMessage
$ws.Range("C2").Value = "Compass Error"
$ws.Range("D2").Value = "0-1"

$ws.Range("B3").Value = $compassMessage
$ws.Range("C3").Value = "Compass disconnected"
$ws.Range("D3").Value = "4-5"

$ws.Range("B4").Value = $compassMessage
$ws.Range("C4").Value = "Compass Error"
$ws.Range("D4").Value = "2-3"

# New rows 5-7: Critical low battery block
$ws.Range("B5").Value = $criticalMessage
$ws.Range("C5").Value = "Critical low battery"
$ws.Range("D5").Value = "0-2"

$ws.Range("B6").Value = $criticalMessage
$ws.Range("C6").Value = "Forced landing in progress"
$ws.Range("D6").Value = "9-12"

$ws.Range("B7").Value = $criticalMessage
$ws.Range("C7").Value = "Aircraft in Auto Power Off Protection"
$ws.Range("D7").Value = "3-8"
